# Update "Fonds de solidarite" volet 1 regional/NAF data for the 2022-06-24 refresh.
# Only the "nombre_aides" (column C) and "montant_total" (column E) figures
# change for a subset of rows; "nombre_entreprises" (column D) stays the same.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = @(
    @{ Row = 10;  NombreAides = 278215; MontantTotal = 1752682990 },
    @{ Row = 107; NombreAides = 26898;  MontantTotal = 36784412 },
    @{ Row = 167; NombreAides = 101542; MontantTotal = 195302952 },
    @{ Row = 168; NombreAides = 285128; MontantTotal = 1213925022 },
    @{ Row = 169; NombreAides = 562684; MontantTotal = 1286541834 },
    @{ Row = 170; NombreAides = 367606; MontantTotal = 2848481747 },
    @{ Row = 171; NombreAides = 115231; MontantTotal = 449035067 },
    @{ Row = 173; NombreAides = 54398;  MontantTotal = 151968108 },
    @{ Row = 174; NombreAides = 357395; MontantTotal = 1020369106 },
    @{ Row = 175; NombreAides = 125779; MontantTotal = 816694666 },
    @{ Row = 177; NombreAides = 96785;  MontantTotal = 174820350 },
    @{ Row = 179; NombreAides = 235818; MontantTotal = 813768825 },
    @{ Row = 180; NombreAides = 141534; MontantTotal = 341253628 },
    @{ Row = 210; NombreAides = 6433;   MontantTotal = 19996314 },
    @{ Row = 266; NombreAides = 71668;  MontantTotal = 219463224 },
    @{ Row = 279; NombreAides = 28969;  MontantTotal = 57090504 },
    @{ Row = 312; NombreAides = 75105;  MontantTotal = 201408498 },
    @{ Row = 313; NombreAides = 220665; MontantTotal = 1371216867 },
    @{ Row = 317; NombreAides = 103596; MontantTotal = 303435659 }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 3).Value = $u.NombreAides
    $ws.Cells.Item($u.Row, 5).Value = $u.MontantTotal
}
